$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WBS")

# Add new category columns to the header row (L1:P1)
$ws.Range("L1").Value = "Schema"
$ws.Range("M1").Value = "Validation"
$ws.Range("N1").Value = "Permissions/Isolation"
$ws.Range("O1").Value = "Workflow"
$ws.Range("P1").Value = "Evidence"

# Rows 2-49 were marked "Done" with a Completed On date (K column).
# Close them out as "Partial" (category-gated done) and clear the
# Completed On date until CI evidence confirms full completion.
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 8).Value = "Partial"
    $ws.Cells.Item($r, 11).ClearContents()
}
